# Apply edits described by the commit: "DOM and Banner author ids added"
# Rows 8-20 were re-permuted (existing records reordered / two records
# relocated) and a cited_by_count value on row 3 was updated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 3
$cell = $ws.Cells.Item(3, 13)
$cell.NumberFormat = "@"
$cell.Value = "9"

# Row 8
$ws.Cells.Item(8, 1).Value = "Ilka Hoof, Klaus Bønnelykke, Thomas Stranzl, Stéphanie Brand, Xingnan Li, Mohamed H. Shamji, Deborah A. Meyers, Eric D. Bateman, Eugene R. Bleecker, Peter S. Andersen"
$ws.Cells.Item(8, 2).Value = "Translational Research, Alk-Abello A/S, Horsholm, Denmark; Copenhagen Prospective Studies on Asthma in Childhood, Copenhagen University Hospital, Copenhagen, Denmark; Translational Research, Alk-Abello A/S, Horsholm, Denmark; Translational Research, Alk-Abello A/S, Horsholm, Denmark; Department of Medicine, The University of Arizona College of Medicine, Tucson, Arizona, USA; National Heart and Lung Institute, Imperial College London, London, UK; Department of Medicine, The University of Arizona College of Medicine, Tucson, Arizona, USA; Division of Respiratory Medicine, Univ of Cape Town, Cape Town, South Africa; Medicine, University of Arizona, Health Sciences Center, Tucson, Arizona, USA; Translational Research, Alk-Abello A/S, Horsholm 2970, Denmark"
$ws.Cells.Item(8, 3).Value = "https://openalex.org/W4390451374"
$ws.Cells.Item(8, 4).Value = "Genetic and T2 biomarkers linked to the efficacy of HDM sublingual immunotherapy in asthma"
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-12-30"
$ws.Cells.Item(8, 6).Value = "Thorax"
$ws.Cells.Item(8, 7).Value = "BMJ"
$ws.Cells.Item(8, 8).Value = "https://doi.org/10.1136/thorax-2023-220707"
$ws.Cells.Item(8, 15).Value = "https://pubmed.ncbi.nlm.nih.gov/38160049"
$ws.Cells.Item(8, 16).Value = "https://doi.org/10.1136/thorax-2023-220707"

# Row 9
$ws.Cells.Item(9, 1).Value = "Ying Zeng, Peiming Zhang, Xingnan Li, Zhan Shi"
$ws.Cells.Item(9, 2).Value = "Guangdong Power Grid Co., Ltd. Power Dispatching and Control Center (China); Guangdong Power Grid Co., Ltd. Power Dispatching and Control Center (China); Guangdong Power Grid Co., Ltd. Power Dispatching and Control Center (China); Guangdong Power Grid Co., Ltd. Power Dispatching and Control Center (China)"
$ws.Cells.Item(9, 3).Value = "https://openalex.org/W4366983696"
$ws.Cells.Item(9, 4).Value = "Service importance-aware virtual embedding strategy for power optical communication network"
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-04-25"
$ws.Cells.Item(9, 6).Value = "N/A"
$ws.Cells.Item(9, 7).Value = "N/A"
$ws.Cells.Item(9, 8).Value = "https://doi.org/10.1117/12.2673169"
$ws.Cells.Item(9, 10).Value = "N/A"
$ws.Cells.Item(9, 11).Value = "closed"
$ws.Cells.Item(9, 15).Value = "NA"
$ws.Cells.Item(9, 16).Value = "https://doi.org/10.1117/12.2673169"

# Row 10
$ws.Cells.Item(10, 1).Value = "Manuel Izquierdo, Chad R. Marion, Frank Genese, John D. Newell, Wanda K. O’Neal, Xingnan Li, Gregory A. Hawkins, Igor Barjaktarević, R. Graham Barr, Stephanie A. Christenson, Christopher B. Cooper, David J. Couper, Jeffrey M. Curtis, MeiLan K. Han, Nadia N. Hansel, Richard E. Kanner, Fernando J. Martínez, Robert Paine, Vickram Tejwani, Prescott G. Woodruff, Joe Zein, Eric A. Hoffman, Stephen P. Peters, Deborah A. Meyers, Eugene R. Bleecker, Victor E. Ortega"
$ws.Cells.Item(10, 2).Value = "; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; "
$ws.Cells.Item(10, 3).Value = "https://openalex.org/W4377008250"
$ws.Cells.Item(10, 4).Value = "Impact of Bronchiectasis on COPD Severity and Alpha-1 Antitrypsin Deficiency as a Risk Factor in Individuals with a Heavy Smoking History"
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-01-01"
$ws.Cells.Item(10, 6).Value = "Chronic obstructive pulmonary diseases"
$ws.Cells.Item(10, 7).Value = "COPD Foundation"
$ws.Cells.Item(10, 8).Value = "https://doi.org/10.15326/jcopdf.2022.0388"
$ws.Cells.Item(10, 10).Value = "publishedVersion"
$ws.Cells.Item(10, 11).Value = "gold"
$ws.Cells.Item(10, 15).Value = "https://pubmed.ncbi.nlm.nih.gov/37199731"
$ws.Cells.Item(10, 16).Value = "https://doi.org/10.15326/jcopdf.2022.0388"

# Row 11
$ws.Cells.Item(11, 1).Value = "Bo Li, Xingnan Li, Prakit Saingam, Tao Yan"
$ws.Cells.Item(11, 2).Value = "Department of Civil and Environmental Engineering, University of Hawaii at Manoa, Honolulu, Hawaii 96822, United States; Department of Civil and Environmental Engineering, University of Hawaii at Manoa, Honolulu, Hawaii 96822, United States; Department of Civil and Environmental Engineering, University of Hawaii at Manoa, Honolulu, Hawaii 96822, United States; Department of Civil and Environmental Engineering, University of Hawaii at Manoa, Honolulu, Hawaii 96822, United States"
$ws.Cells.Item(11, 3).Value = "https://openalex.org/W4383533147"
$ws.Cells.Item(11, 4).Value = "Understanding the Microbiological Quality of Drinking Water at the Point of Consumption with Citizen Science"
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-07-07"
$ws.Cells.Item(11, 6).Value = "ACS ES&T Water"
$ws.Cells.Item(11, 7).Value = "American Chemical Society"
$ws.Cells.Item(11, 8).Value = "https://doi.org/10.1021/acsestwater.3c00184"
$ws.Cells.Item(11, 10).Value = "N/A"
$ws.Cells.Item(11, 11).Value = "closed"
$ws.Cells.Item(11, 15).Value = "NA"
$ws.Cells.Item(11, 16).Value = "https://doi.org/10.1021/acsestwater.3c00184"

# Row 12
$ws.Cells.Item(12, 1).Value = "Xiaozhi Deng, Bo Li, Xingnan Li, Zhentian Wu, Zhihua Yang"
$ws.Cells.Item(12, 2).Value = "Guangdong Power Grid Co., Ltd.,Power Dispatching and Control Center,Communication Management Department,Guangdong,China; Guangdong Power Grid Co., Ltd.,Power Dispatching and Control Center,Communication Management Department,Guangdong,China; Guangdong Power Grid Co., Ltd.,Power Dispatching and Control Center,Communication Management Department,Guangdong,China; Guangdong Electric Power Communication Technology Co. Ltd.,Science Innovation Department,Guangdong,China; Guangdong Electric Power Communication Technology Co. Ltd.,Science Innovation Department,Guangdong,China"
$ws.Cells.Item(12, 3).Value = "https://openalex.org/W4383988695"
$ws.Cells.Item(12, 4).Value = "Container and Microservice-Based Resource Management for Distribution Station Area"
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-05-19"
$ws.Cells.Item(12, 6).Value = "N/A"
$ws.Cells.Item(12, 7).Value = "N/A"
$ws.Cells.Item(12, 8).Value = "https://doi.org/10.1109/icmsp58539.2023.10170832"
$ws.Cells.Item(12, 15).Value = "NA"
$ws.Cells.Item(12, 16).Value = "https://doi.org/10.1109/icmsp58539.2023.10170832"

# Row 13
$ws.Cells.Item(13, 1).Value = "Xingnan Li, Xiaozhi Deng, Zhan Shi, Zhihua Yang, Xin Qian"
$ws.Cells.Item(13, 2).Value = "Communication Management Department, Power Dispatching Control Center, Guangdong Power Grid Co., Ltd., Guangdong, China; Communication Management Department, Power Dispatching Control Center, Guangdong Power Grid Co., Ltd., Guangdong, China; Communication Management Department, Power Dispatching Control Center, Guangdong Power Grid Co., Ltd., Guangdong, China; Science Innovation Department, Guangdong Electric Power Communication Technology Co. Ltd., Guangdong, China; Science Innovation Department, Guangdong Electric Power Communication Technology Co. Ltd., Guangdong, China"
$ws.Cells.Item(13, 3).Value = "https://openalex.org/W4384026380"
$ws.Cells.Item(13, 4).Value = "Reliable Aggregation Method of Power Line Communication Subcarriers in Complex Power Electronic Colored Noise Environment"
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-05-19"
$ws.Cells.Item(13, 6).Value = "N/A"
$ws.Cells.Item(13, 7).Value = "N/A"
$ws.Cells.Item(13, 8).Value = "https://doi.org/10.1109/icmsp58539.2023.10170963"
$ws.Cells.Item(13, 10).Value = "N/A"
$ws.Cells.Item(13, 11).Value = "closed"
$ws.Cells.Item(13, 16).Value = "https://doi.org/10.1109/icmsp58539.2023.10170963"

# Row 14
$ws.Cells.Item(14, 1).Value = "Manuel Izquierdo, Chad R. Marion, Frank Genese, John D. Newell, Wanda K. O’Neal, Xingnan Li, Gregory A. Hawkins, Igor Barjaktarević, R. Graham Barr, Stephanie A. Christenson, Christopher B. Cooper, David J. Couper, Jeffrey M. Curtis, MeiLan K. Han, Nadia N. Hansel, Richard E. Kanner, Fernando J. Martínez, Robert Paine, Vickram Tejwani, Prescott G. Woodruff, Joe Zein, Eric A. Hoffman, Stephen P. Peters, Deborah A. Meyers, Eugene R. Bleecker, Victor E. Ortega"
$ws.Cells.Item(14, 2).Value = "; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; "
$ws.Cells.Item(14, 3).Value = "https://openalex.org/W4385267546"
$ws.Cells.Item(14, 4).Value = "Impact of Bronchiectasis on COPD Severity and Alpha-1 Antitrypsin Deficiency as a Risk Factor in Individuals with a Heavy Smoking History"
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-01-01"
$ws.Cells.Item(14, 6).Value = "Chronic obstructive pulmonary diseases"
$ws.Cells.Item(14, 7).Value = "COPD Foundation"
$ws.Cells.Item(14, 8).Value = "https://doi.org/10.15326/jcopdf.2023.0388"
$ws.Cells.Item(14, 10).Value = "publishedVersion"
$ws.Cells.Item(14, 11).Value = "gold"
$ws.Cells.Item(14, 15).Value = "https://pubmed.ncbi.nlm.nih.gov/37199731"
$ws.Cells.Item(14, 16).Value = "https://doi.org/10.15326/jcopdf.2023.0388"

# Row 15
$ws.Cells.Item(15, 1).Value = "Huashi Li, Xingnan Li"
$ws.Cells.Item(15, 2).Value = "Statistics Consulting Lab, BIO5 Institute, University of Arizona, Tucson, AZ, USA; Division of Genetics, Genomics and Precision Medicine, Department of Medicine, University of Arizona College of Medicine, Tucson, AZ, USA"
$ws.Cells.Item(15, 3).Value = "https://openalex.org/W4385716105"
$ws.Cells.Item(15, 4).Value = "Genetic relationships between high blood eosinophil count, asthma susceptibility, and asthma severity"
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-08-22"
$ws.Cells.Item(15, 6).Value = "Journal of Asthma"
$ws.Cells.Item(15, 7).Value = "Taylor & Francis"
$ws.Cells.Item(15, 8).Value = "https://doi.org/10.1080/02770903.2023.2247490"
$ws.Cells.Item(15, 15).Value = "https://pubmed.ncbi.nlm.nih.gov/37560908"
$ws.Cells.Item(15, 16).Value = "https://doi.org/10.1080/02770903.2023.2247490"

# Row 16
$ws.Cells.Item(16, 1).Value = "Zhan Shi, Xiaozhi Deng, Xingnan Li, Zhihua Yang, Xiaohong Qian"
$ws.Cells.Item(16, 2).Value = "Communication Management Department, Guangdong Power Grid Co., Ltd. Power Dispatching Control Center, Guangdong, China; Communication Management Department, Guangdong Power Grid Co., Ltd. Power Dispatching Control Center, Guangdong, China; Communication Management Department, Guangdong Power Grid Co., Ltd. Power Dispatching Control Center, Guangdong, China; ; "
$ws.Cells.Item(16, 3).Value = "https://openalex.org/W4387049368"
$ws.Cells.Item(16, 4).Value = "Value Gain-based Power Line Subcarrier Aggregation Method under Colored Noise"
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-09-01"
$ws.Cells.Item(16, 6).Value = "Journal of Physics: Conference Series"
$ws.Cells.Item(16, 7).Value = "IOP Publishing"
$ws.Cells.Item(16, 8).Value = "https://doi.org/10.1088/1742-6596/2584/1/012132"
$ws.Cells.Item(16, 9).Value = "N/A"
$ws.Cells.Item(16, 11).Value = "bronze"
$ws.Cells.Item(16, 15).Value = "NA"
$ws.Cells.Item(16, 16).Value = "https://doi.org/10.1088/1742-6596/2584/1/012132"

# Row 17
$ws.Cells.Item(17, 1).Value = "Xingnan Li, Jiangang Lü, Peiming Zhang"
$ws.Cells.Item(17, 2).Value = "Power Dispatching Control Center of Guangdong Power Grid Co., Ltd,Guangzhou,China; Power Dispatching Control Center of Guangdong Power Grid Co., Ltd,Guangzhou,China; Power Dispatching Control Center of Guangdong Power Grid Co., Ltd,Guangzhou,China"
$ws.Cells.Item(17, 3).Value = "https://openalex.org/W4387251115"
$ws.Cells.Item(17, 4).Value = "A GNN-Based Routing and Scheduling Mechanism for Multi-domain Computing First Network"
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-08-12"
$ws.Cells.Item(17, 8).Value = "https://doi.org/10.1109/ccis59572.2023.10262822"
$ws.Cells.Item(17, 16).Value = "https://doi.org/10.1109/ccis59572.2023.10262822"

# Row 18
$ws.Cells.Item(18, 1).Value = "Natalie Iannuzo, Holly Welfley, Nicholas Li, Michael G. Johnson, Stefano Guerra, Xingnan Li, Darren A. Cusanovich, Paul Langlais, Julie G. Ledford"
$ws.Cells.Item(18, 2).Value = "University of Arizona, Tucson, United States; University of Arizona, Tucson, United States; Basis Tucson North, Tucson, United States; University of Arizona, Tucson, United States; University of Arizona, Tucson, United States; University of Arizona, Tucson, United States; University of Arizona, Tucson, United States; University of Arizona, Tucson, United States; University of Arizona, Tucson, United States"
$ws.Cells.Item(18, 3).Value = "https://openalex.org/W4387979514"
$ws.Cells.Item(18, 4).Value = "CC16 Induces Pulmonary Epithelial-Driven SPLUNC1 Expression by Signaling through VLA-2."
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-09-09"
$ws.Cells.Item(18, 6).Value = "N/A"
$ws.Cells.Item(18, 7).Value = "N/A"
$ws.Cells.Item(18, 8).Value = "https://doi.org/10.1183/13993003.congress-2023.oa4964"
$ws.Cells.Item(18, 15).Value = "NA"
$ws.Cells.Item(18, 16).Value = "https://doi.org/10.1183/13993003.congress-2023.oa4964"

# Row 19
$ws.Cells.Item(19, 1).Value = "Natalie Iannuzo, Holly Welfley, Nicholas Li, Michael D. L. Johnson, Joselyn Rojas, Francesca Polverino, Stefano Guerra, Xingnan Li, Darren A. Cusanovich, Paul Langlais, Julie G. Ledford"
$ws.Cells.Item(19, 2).Value = "Department of Cellular and Molecular Medicine, University of Arizona, Tucson, AZ, United States; Asthma and Airway Disease Research Center, Tucson, AZ, United States; BASIS Tucson North, Tucson, AZ, United States; Department of Immunobiology, University of Arizona, Tucson, AZ, United States; Baylor College of Medicine, Houston, TX, United States; Baylor College of Medicine, Houston, TX, United States; Asthma and Airway Disease Research Center, Tucson, AZ, United States; Department of Medicine, Division of Pulmonary, Allergy, Critical Care, and Sleep Medicine, University of Arizona, Tucson, AZ, United States; Department of Medicine, Division of Genetics, Genomics, and Precision Medicine, University of Arizona, Tucson, AZ, United States; Asthma and Airway Disease Research Center, Tucson, AZ, United States; Department of Cellular and Molecular Medicine, University of Arizona, Tucson, AZ, United States; Department of Medicine, Division of Endocrinology, University of Arizona, Tucson, AZ, United States; Asthma and Airway Disease Research Center, Tucson, AZ, United States; Department of Cellular and Molecular Medicine, University of Arizona, Tucson, AZ, United States"
$ws.Cells.Item(19, 3).Value = "https://openalex.org/W4388832621"
$ws.Cells.Item(19, 4).Value = "CC16 drives VLA-2-dependent SPLUNC1 expression"
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-11-20"
$ws.Cells.Item(19, 6).Value = "Frontiers in Immunology"
$ws.Cells.Item(19, 7).Value = "Frontiers Media"
$ws.Cells.Item(19, 8).Value = "https://doi.org/10.3389/fimmu.2023.1277582"
$ws.Cells.Item(19, 9).Value = "cc-by"
$ws.Cells.Item(19, 10).Value = "publishedVersion"
$ws.Cells.Item(19, 11).Value = "gold"
$ws.Cells.Item(19, 15).Value = "https://pubmed.ncbi.nlm.nih.gov/38053993"
$ws.Cells.Item(19, 16).Value = "https://doi.org/10.3389/fimmu.2023.1277582"

# Row 20
$ws.Cells.Item(20, 1).Value = "Peiming Zhang, Xingnan Li, Yuanjie Liu"
$ws.Cells.Item(20, 2).Value = "Guangdong Power Grid Co.,Ltd,Power Dispatching Control Center,Guangdong,China; Guangdong Power Grid Co.,Ltd,Power Dispatching Control Center,Guangdong,China; Guangdong Power Grid Co.,Ltd,Power Dispatching Control Center,Guangdong,China"
$ws.Cells.Item(20, 3).Value = "https://openalex.org/W4390187821"
$ws.Cells.Item(20, 4).Value = "Routing Optimization Mechanism for SRv6 Based Power Data Network"
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = "2023-09-27"
$ws.Cells.Item(20, 8).Value = "https://doi.org/10.1109/wsce59557.2023.10365774"
$ws.Cells.Item(20, 16).Value = "https://doi.org/10.1109/wsce59557.2023.10365774"
